$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff shows that the weekly data update:
#  - appends a new row (58) with the same data that row 57 currently has
#  - updates row 57 in place with this week's new values (date, quality, volume)
#
# So first, copy the current (old) row 57 values down into the new row 58,
# preserving the date number format on column D, then overwrite row 57
# with the new week's figures.

$lastCol = 18  # columns A..R

for ($col = 1; $col -le $lastCol; $col++) {
    $src = $ws.Cells.Item(57, $col)
    $dst = $ws.Cells.Item(58, $col)
    $dst.Value = $src.Value2
}
$ws.Cells.Item(58, 4).NumberFormat = $ws.Cells.Item(57, 4).NumberFormat

# Now update row 57 with the new week's values
$ws.Cells.Item(57, 4).Value = 44628   # Fecha
$ws.Cells.Item(57, 9).Value = "Segunda"  # Calidad
$ws.Cells.Item(57, 10).Value = 200    # Volumen

$ws.Range("A1").Select()
